$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row below the current row 18 ("مبرد قدم"), shifting the
#    totals row (19->20) and the footer row (20->21) down by one.
$ws.Rows("19:19").Insert()

# 2) Clone row 18's formatting into the freshly inserted row 19 (reuses the
#    existing style records instead of minting new ones).
$ws.Range("A18:N18").Copy()
$ws.Range("A19:N19").PasteSpecial(-4122)

# 3) Clone row 18's values into row 19 too - this duplicates the old
#    "مبرد قدم" product line into the new row.
$ws.Range("A18:N18").Copy()
$ws.Range("A19:N19").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# 4) Re-create the three merged ranges for the new row, matching the pattern
#    used by every other product row.
$ws.Range("B19:G19").Merge()
$ws.Range("H19:K19").Merge()
$ws.Range("L19:M19").Merge()
$ws.Rows("19:19").RowHeight = 24.75

# 5) Row 19 is product #16 ("مبرد قدم") - fix up its sequence number.
$ws.Range("A19").Value2 = 16

# 6) Row 18 becomes the brand-new product "لزقه النمر بسعر القطعه" (#15).
$ws.Range("B18").Value2 = "لزقه النمر بسعر القطعه"
$ws.Range("H18").Value2 = "50:0"
$ws.Range("L18").Value2 = -15
$ws.Range("N18").Value2 = "1:0"

# 7) Update the recalculated total (old row 19, now row 20) and its height.
$ws.Range("K20").Value2 = 329.5
$ws.Rows("20:20").RowHeight = 26.25
